# V2 started, using different crystal, ICM20948, 10 pin SWD header, and
# multiple component swaps. Appends new requirement/notes rows to the
# "Sheet1" requirements list (col A), leaving a blank row between each
# logical group, matching the order the notes were typed in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "all bottom components..." note (row 25, rows 26-27 left blank) ---
$ws.Range("A25").Value = "all bottom components easy to solder, simple"

# --- IMU / crystal swap notes (rows 28-31, rows 32-33 left blank) ---
$ws.Range("A28").Value = "as many as possible are simple"
$ws.Range("A29").Value = "using ICM20948 because MPU6050 has no SPI bus"
$ws.Range("A30").Value = "level shifter IC in use is merely N-channel MOSFET with sufficiently low Gate threshold voltage"
$ws.Range("A31").Value = "switched crystal to 16MHz oscillator referenced in video"

# --- "parts not in database at all" block (rows 34-36) ---
$ws.Range("A35").Value = "ICM20948"
$ws.Range("A34").Value = "parts not in database at all:"
$ws.Range("A36").Value = "tactile reset switch (unless we use an extended part)"

# --- "parts not basic" block, interleaved as it was authored (rows 37-46, rows 41-42 left blank) ---
$ws.Range("A43").Value = "parts not basic:"
$ws.Range("A44").Value = "STM32F446RET6"

$ws.Range("A37").Value = "SD card connector"
$ws.Range("A38").Value = "USB connector"

$ws.Range("A45").Value = "FT230XQ"
$ws.Range("A46").Value = "1.8V regulator"

$ws.Range("A39").Value = "motor control MOSFET"
$ws.Range("A40").Value = "all connectors"

# --- best-effort viewport scroll to match the author's saved window state ---
# (topLeftCell="A19"); harmless no-op if the host doesn't persist it.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

# --- final selection, matches the author's saved cursor position ---
$ws.Range("D35").Select() | Out-Null
